# Update on cleaning data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows appended at the bottom of the table (rows 402-433) ---
$aids  = @(501,502,503,504,505,506,507,508,509,510,511,512,513,514,515,516,517,518,519,520,521,522,523,524,525,526,527,528,529,530,551,561)
$evals = @(2.1,2.2,2.1,2.2,2.2,2.0,2.1,2.2,2.2,2.2,2.2,2.3,2.1,2.1,2.3,2.1,2.3,2.1,2.3,2.0,2.5,2.2,2.3,2.1,2.1,2.3,2.2,2.2,2.2,2.2,2.2,2.2)

$startRow = 402
for ($i = 0; $i -lt $aids.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $aids[$i]
    $ws.Cells.Item($r, 5).Value = $evals[$i]
}

# --- New note marker on row 20 (column G), a backtick character ---
$ws.Range("G20").Value = "``"

# --- Update the selected range shown when the file is opened ---
$ws.Range("A163:XFD163").Select()
